$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend Table1 with three new columns (auto-extends the table range) ---
$lo = $ws.ListObjects.Item(1)
$lo.ListColumns.Add() | Out-Null
$lo.ListColumns.Add() | Out-Null
$lo.ListColumns.Add() | Out-Null

# --- Correct the "Total Cycles To Execute" value for the hardware-scheduled
#     pipeline / Grendle row ---
$ws.Cells.Item(6, 4).Value = 4459

# --- Copy existing formatting onto the new columns so the new columns match
#     the look of the rest of the table ---
$ws.Range("C1").Copy() | Out-Null
$ws.Range("E1:G1").PasteSpecial(-4122) | Out-Null

$ws.Range("D2").Copy() | Out-Null
$ws.Range("E2:G2").PasteSpecial(-4122) | Out-Null

$ws.Range("D3").Copy() | Out-Null
$ws.Range("E3:G3").PasteSpecial(-4122) | Out-Null

$ws.Range("D4").Copy() | Out-Null
$ws.Range("E4:G4").PasteSpecial(-4122) | Out-Null

$ws.Range("D5").Copy() | Out-Null
$ws.Range("E5:G5").PasteSpecial(-4122) | Out-Null

$ws.Range("D6").Copy() | Out-Null
$ws.Range("E6:G6").PasteSpecial(-4122) | Out-Null

$ws.Range("D7").Copy() | Out-Null
$ws.Range("E7:G7").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Store the new column data as text so values such as "1.00" keep their
#     trailing zeros ---
$ws.Range("E2:G7").NumberFormat = "@"

# --- Header row ---
$ws.Range("E1").Value = "CPI"
$ws.Range("F1").Value = "Max Cycle Time"
$ws.Range("G1").Value = "Total Execution time"

# --- CPI column ---
$ws.Range("E2").Value = "1.00"
$ws.Range("E3").Value = "1.00"
$ws.Range("E4").Value = "1.15"
$ws.Range("E5").Value = "1.26"
$ws.Range("E6").Value = "2.11"
$ws.Range("E7").Value = "2.26"

# --- Max Cycle Time column (only populated on the first row of each
#     processor/benchmark pair) ---
$ws.Range("F2").Value = "24.74mhz"
$ws.Range("F4").Value = "53.02mhz"
$ws.Range("F6").Value = "47.39mhz"

# --- Total Execution time column ---
$ws.Range("G2").Value = "85,529.5ns"
$ws.Range("G3").Value = "33,427.6ns"
$ws.Range("G4").Value = "114,843ns"
$ws.Range("G5").Value = "45,662ns"
$ws.Range("G6").Value = "94091.6ns"
$ws.Range("G7").Value = "39375.4ns"

# --- Column widths ---
$ws.Range("C1").EntireColumn.ColumnWidth = 16.6640625
$ws.Range("D1").EntireColumn.ColumnWidth = 22.6640625
$ws.Range("E1").EntireColumn.ColumnWidth = 8.796875
$ws.Range("F1").EntireColumn.ColumnWidth = 18.3984375
$ws.Range("G1").EntireColumn.ColumnWidth = 22.6640625

# --- Restore the cosmetic selection state ---
$ws.Range("F13").Select() | Out-Null

$wb.Save()
Write-Output "edit.ps1 completed"
